$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.605.92'
$ws.Range("E2").Value = '  +0.89%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.769.28'
$ws.Range("E3").Value = '  -1.09%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.25%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '442.66'
$ws.Range("E5").Value = '  +5.15%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.38'
$ws.Range("E6").Value = '  +11.52%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.617'
$ws.Range("E7").Value = '  +2.79%  '

$ws.Range("E8").Value = '  +0.18%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.727'
$ws.Range("E9").Value = '  +1.95%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.149'
$ws.Range("E10").Value = '  -8.18%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000307'
$ws.Range("E11").Value = '  -10.71%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '42.60'
$ws.Range("E12").Value = '  +6.85%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.28'
$ws.Range("E13").Value = '  +4.31%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.368.11'
$ws.Range("E14").Value = '  -0.71%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.74'
$ws.Range("E15").Value = '  -5.68%  '

$ws.Range("E16").Value = '  -0.26%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.786.07'
$ws.Range("E17").Value = '  -1.00%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '19.73'
$ws.Range("E18").Value = '  +1.96%  '

$ws.Range("E19").Value = '  +6.88%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '66.630.94'
$ws.Range("E20").Value = '  +0.78%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '409.62'
$ws.Range("E21").Value = '  +2.48%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.41'
$ws.Range("E22").Value = '  +1.68%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.24'
$ws.Range("E23").Value = '  +9.41%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.94'
$ws.Range("E24").Value = '  +1.88%  '

$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.38'
$ws.Range("E25").Value = '  +6.58%  '

$ws.Range("B26").Value = 'EthereumClassic'
$ws.Range("C26").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '36.60'
$ws.Range("E26").Value = '  -0.47%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '5.59'
$ws.Range("E27").Value = '  -2.71%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.71'
$ws.Range("E28").Value = '  +33.33%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.62'
$ws.Range("E29").Value = '  +3.40%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '730.86'
$ws.Range("E30").Value = '  +5.64%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.84'
$ws.Range("E31").Value = '  +13.78%  '

$ws.Range("E32").Value = '  +11.74%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.73'
$ws.Range("E33").Value = '  -0.03%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '43.21'
$ws.Range("E34").Value = '  +15.22%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.156'
$ws.Range("E35").Value = '  +4.61%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '56.31'
$ws.Range("E36").Value = '  +3.10%  '

$ws.Range("B37").Value = 'Dai'
$ws.Range("C37").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.998'
$ws.Range("E37").Value = '  -0.19%  '

$ws.Range("B38").Value = 'NEARProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.51'
$ws.Range("E38").Value = '  +25.01%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0472'
$ws.Range("E39").Value = '  +5.70%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.64'
$ws.Range("E40").Value = '  +33.87%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.87'
$ws.Range("E41").Value = '  -1.03%  '

$ws.Range("E42").Value = '  +0.37%  '

$ws.Range("B43").Value = 'Stellar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.139'
$ws.Range("E43").Value = '  +4.27%  '

$ws.Range("B44").Value = 'ApeXProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.30'
$ws.Range("E44").Value = '  +7.63%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.333'
$ws.Range("E45").Value = '  +17.24%  '

$ws.Range("E46").Value = '  -12.20%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.31'
$ws.Range("E47").Value = '  +1.95%  '

$ws.Range("B48").Value = 'WEMIXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.65'
$ws.Range("E48").Value = '  +5.08%  '

$ws.Range("B49").Value = 'ARBITRUM'
$ws.Range("C49").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.07'
$ws.Range("E49").Value = '  +2.03%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '142.51'
$ws.Range("E50").Value = '  -1.05%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.80'
$ws.Range("E51").Value = '  +3.47%  '
